$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.335.14"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +15.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.686.10"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +9.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -1.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.91"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +9.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9975"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3729"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +3.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3437"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +8.46%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.35"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +18.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.186"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +8.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07310"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +7.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.62"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +10.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.105"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +8.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.762"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +6.64%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.673.20"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +9.59%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001109"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +6.37%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9980"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.36%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06726"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +10.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "81.53"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +12.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.49"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +10.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.131"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +7.75%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.05"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.222.05"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +14.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.406"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.683"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +21.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.359"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -9.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.31"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +2.79%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.56"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +11.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.858.24"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +9.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.33"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +7.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.474"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +25.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.038"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.46%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9918"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +16.64%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.744"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +15.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08463"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.92%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "12.46"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +17.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.381"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +8.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06428"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +9.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.888"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +15.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.291"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.02339"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +11.51%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.2114"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +10.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6157"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +13.49%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9970"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.801"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +6.22%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "13.20"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.56%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5967"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +9.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "127.96"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +5.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.024"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +8.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07170"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +9.22%  "
